# Apply the commit: "update spreadsheet, modify program to read data from
# multiple worksheets" — bumps the raw input figures on the "2017" sheet
# from hundreds to thousands (the dependent SUM() formulas recalc on
# their own) and flips which sheet/range is active & selected.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) New leaf (non-formula) input values for rows 5-9 on sheet "2017".
#    Every one of the five data rows uses the same column -> new-value
#    mapping, so loop rows x columns.
# ---------------------------------------------------------------------
$ws2017 = $wb.Worksheets.Item("2017")

$newValues = [ordered]@{
    "F"  = 3000;  "G"  = 4000
    "I"  = 5000;  "J"  = 6000
    "L"  = 7000;  "M"  = 8000
    "R"  = 9000;  "S"  = 10000
    "U"  = 11000; "V"  = 12000
    "X"  = 13000; "Y"  = 14000
    "AD" = 15000; "AE" = 16000
    "AG" = 17000; "AH" = 18000
    "AJ" = 19000; "AK" = 20000
    "AP" = 21000; "AQ" = 22000
    "AS" = 23000; "AT" = 24000
    "AV" = 25000; "AW" = 26000
}

foreach ($row in 5..9) {
    foreach ($col in $newValues.Keys) {
        $ws2017.Range("$col$row").Value = $newValues[$col]
    }
}

# ---------------------------------------------------------------------
# 2) View/selection state: "2017" becomes the active sheet/tab (was
#    "2016"), scrolled toward column AA with AV5:AW9 selected; "2016"
#    keeps its selection anchored at AW9 once it's no longer active.
# ---------------------------------------------------------------------
$ws2016 = $wb.Worksheets.Item("2016")
$ws2016.Range("AW9").Select()
$excel.ActiveWindow.ScrollColumn = 29

$ws2017.Activate()
$excel.ActiveWindow.ScrollColumn = 27
$excel.ActiveWindow.ScrollRow = 1
$ws2017.Range("AV5:AW9").Select()
